# Applies numeric corrections to various crafting-profit rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 38799.6
$ws.Range("I54").Value = 21999
$ws.Range("J54").Value = 50000
$ws.Range("K54").Value = 21999
$ws.Range("L54").Value = 50000
$ws.Range("M54").Value = -21513
$ws.Range("N54").Value = -50972

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1623
$ws.Range("I107").Value = 1920.1428
$ws.Range("J107").Value = 583
$ws.Range("K107").Value = 1920.1428
$ws.Range("L107").Value = 583
$ws.Range("M107").Value = -0.1428000000000793
$ws.Range("N107").Value = -4423

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2421.9443
$ws.Range("I137").Value = 1842.2858
$ws.Range("J137").Value = 4450.75
$ws.Range("K137").Value = 5526.857400000001
$ws.Range("L137").Value = 13352.25
$ws.Range("M137").Value = -2976.857400000001
$ws.Range("N137").Value = -18452.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7268.5615
$ws.Range("I32").Value = 6064.4614
$ws.Range("J32").Value = 9877.444
$ws.Range("K32").Value = 6064.4614
$ws.Range("L32").Value = 9877.444
$ws.Range("M32").Value = -5777.4614
$ws.Range("N32").Value = -10451.444

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5543495.5
$ws.Range("I63").Value = 8659218
$ws.Range("J63").Value = 4433.3335
$ws.Range("K63").Value = 8659218
$ws.Range("L63").Value = 4433.3335
$ws.Range("M63").Value = -8658532
$ws.Range("N63").Value = -5805.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5543495.5
$ws.Range("I66").Value = 8659218
$ws.Range("J66").Value = 4433.3335
$ws.Range("K66").Value = 43296090
$ws.Range("L66").Value = 22166.6675
$ws.Range("M66").Value = -43292658
$ws.Range("N66").Value = -29030.6675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 964.2
$ws.Range("I110").Value = 964.2
$ws.Range("K110").Value = 964.2
$ws.Range("M110").Value = 1080.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2410.4468
$ws.Range("I132").Value = 1878.6857
$ws.Range("J132").Value = 3961.4167
$ws.Range("K132").Value = 5636.0571
$ws.Range("L132").Value = 11884.2501
$ws.Range("M132").Value = -3106.0571
$ws.Range("N132").Value = -16944.2501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10641620
$ws.Range("I31").Value = 2087.9656
$ws.Range("J31").Value = 27783088
$ws.Range("K31").Value = 2087.9656
$ws.Range("L31").Value = 27783088
$ws.Range("M31").Value = -1792.9656
$ws.Range("N31").Value = -27783678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10641620
$ws.Range("I34").Value = 2087.9656
$ws.Range("J34").Value = 27783088
$ws.Range("K34").Value = 2087.9656
$ws.Range("L34").Value = 27783088
$ws.Range("M34").Value = -1885.9656
$ws.Range("N34").Value = -27783492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15388949
$ws.Range("I99").Value = 25002166
$ws.Range("J99").Value = 7802.8
$ws.Range("K99").Value = 25002166
$ws.Range("L99").Value = 7802.8
$ws.Range("M99").Value = -25000668
$ws.Range("N99").Value = -10798.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 15388949
$ws.Range("I126").Value = 25002166
$ws.Range("J126").Value = 7802.8
$ws.Range("K126").Value = 75006498
$ws.Range("L126").Value = 23408.4
$ws.Range("M126").Value = -75004028
$ws.Range("N126").Value = -28348.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1609.2106
$ws.Range("I5").Value = 492.8889
$ws.Range("J5").Value = 2613.9
$ws.Range("K5").Value = 1478.6667
$ws.Range("L5").Value = 7841.700000000001
$ws.Range("M5").Value = -1366.6667
$ws.Range("N5").Value = -8065.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 2775
$ws.Range("J101").Value = 2775
$ws.Range("L101").Value = 8325
$ws.Range("N101").Value = -13193

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 72104.86
$ws.Range("I107").Value = 592.6
$ws.Range("K107").Value = 1777.8
$ws.Range("M107").Value = 142.1999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 618.82355
$ws.Range("I113").Value = 552
$ws.Range("J113").Value = 678.2222
$ws.Range("K113").Value = 1656
$ws.Range("L113").Value = 2034.6666
$ws.Range("M113").Value = 514
$ws.Range("N113").Value = -6374.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1576.1724
$ws.Range("I132").Value = 783.375
$ws.Range("K132").Value = 7050.375
$ws.Range("M132").Value = -4520.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1609.2106
$ws.Range("I135").Value = 492.8889
$ws.Range("J135").Value = 2613.9
$ws.Range("K135").Value = 4436.0001
$ws.Range("L135").Value = 23525.1
$ws.Range("M135").Value = -1901.0001
$ws.Range("N135").Value = -28595.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 103
$ws.Range("I2").Value = 83.333336
$ws.Range("J2").Value = 162
$ws.Range("K2").Value = 83.333336
$ws.Range("L2").Value = 162
$ws.Range("M2").Value = 29.666664
$ws.Range("N2").Value = -388

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1667.1714
$ws.Range("I102").Value = 1173.9584
$ws.Range("J102").Value = 2743.2727
$ws.Range("K102").Value = 1173.9584
$ws.Range("L102").Value = 2743.2727
$ws.Range("M102").Value = 448.0416
$ws.Range("N102").Value = -5987.2727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 31500
$ws.Range("J105").Value = 31500
$ws.Range("L105").Value = 31500
$ws.Range("N105").Value = -38488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2982
$ws.Range("I122").Value = 1413.375
$ws.Range("K122").Value = 4240.125
$ws.Range("M122").Value = -1790.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 956.5714
$ws.Range("I55").Value = 424.5
$ws.Range("J55").Value = 1666
$ws.Range("K55").Value = 424.5
$ws.Range("L55").Value = 1666
$ws.Range("M55").Value = -251.5
$ws.Range("N55").Value = -2012

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1219.2941
$ws.Range("I61").Value = 1034.8667
$ws.Range("J61").Value = 2602.5
$ws.Range("K61").Value = 1034.8667
$ws.Range("L61").Value = 2602.5
$ws.Range("M61").Value = -832.8667
$ws.Range("N61").Value = -3006.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4446469.5
$ws.Range("I93").Value = 6945871.5
$ws.Range("J93").Value = 3088.6667
$ws.Range("K93").Value = 6945871.5
$ws.Range("L93").Value = 3088.6667
$ws.Range("M93").Value = -6944623.5
$ws.Range("N93").Value = -5584.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1219.2941
$ws.Range("I113").Value = 1034.8667
$ws.Range("J113").Value = 2602.5
$ws.Range("K113").Value = 1034.8667
$ws.Range("L113").Value = 2602.5
$ws.Range("M113").Value = 1135.1333
$ws.Range("N113").Value = -6942.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 3342
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 483.8095
$ws.Range("I100").Value = 490.3889
$ws.Range("J100").Value = 444.33334
$ws.Range("K100").Value = 980.7778
$ws.Range("L100").Value = 888.66668
$ws.Range("M100").Value = -439.7778
$ws.Range("N100").Value = -1970.66668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 631.1429000000001
$ws.Range("I107").Value = 703
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 2109
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = -189
$ws.Range("N107").Value = -4440
